# Updates cryptocurrency price/volume data to match the latest scrape.
# Columns: B=Coin name, C=Link, D=Price, E=Volume(1h) change %.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '29.948.53'
$ws.Range('E2').Value = '  +0.31%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.893.69'
$ws.Range('E3').Value = '  -0.07%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.01%  '

# Row 5: XRP
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7722'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.37%  '

# Row 6: BNB
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.37%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.04%  '

# Row 8: Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3128'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.80%  '

# Row 9: Solana
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.77'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.76%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07270'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.58%  '

# Row 11: TRON
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08044'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.65%  '

# Row 12: Polygon
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7711'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.39%  '

# Row 13: Polkadot
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.461'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.90%  '

# Row 14: WrappedEther
$ws.Range('D14').Value = '1.916.77'
$ws.Range('E14').Value = '  +1.26%  '

# Row 15: Litecoin
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '95.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.80%  '

# Row 16: Uniswap
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.182'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.22%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '29.934.60'
$ws.Range('E17').Value = '  +0.29%  '

# Row 18: Avalanche
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.74%  '

# Row 19: BitcoinCash
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '

# Row 20: ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007858'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.95%  '

# Row 21: Chainlink
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.168.04'
$ws.Range('E21').Value = '  +1.83%  '

# Row 22: Dai
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '

# Row 23: WrappedliquidstakedEther2.0
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.115'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.17%  '

# Row 24: BinanceUSD
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.01%  '

# Row 25: Stellar
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1583'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.08%  '

# Row 26: Cosmos
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.526'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.21%  '

# Row 27: Monero
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.33%  '

# Row 28: EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.39%  '

# Row 29: LidoDAOToken
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.038'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.87%  '

# Row 30: Toncoin
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.415'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.81%  '

# Row 31: PancakeSwap
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.543'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.19%  '

# Row 32: Filecoin
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.519'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.08%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.093'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.73%  '

# Row 34: Hedera
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05496'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.35%  '

# Row 35: ARBITRUM
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.245'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.82%  '

# Row 36: ImmutableX
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7488'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.98%  '

# Row 37: Frax
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.003'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.34%  '

# Row 38: HuobiToken
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.689'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.96%  '

# Row 39: VeChain
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01934'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.59%  '

# Row 40: MXToken
$ws.Range('E40').Value = '  +0.27%  '

# Row 41: TheSandbox
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4491'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.72%  '

# Row 42: Aave
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.62%  '

# Row 43: FraxShare
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.099.95'
$ws.Range('E43').Value = '  -1.89%  '

# Row 44: Maker
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.073'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.24%  '

# Row 45: TrustWalletToken
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8520'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.27%  '

# Row 46: PaxDollar
$ws.Range('E46').Value = '  +0.03%  '

# Row 47: RenderToken
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.892'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.68%  '

# Row 48: Quant
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.85%  '

# Row 49: Aptos
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.598'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.08%  '

# Row 50: EnergySwap
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.799'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.59%  '

# Row 51: SynthetixNetwork
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.011'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.08%  '
